$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing contingency rows (row r: C=from_bus-ish val, D=to_bus-ish val, E=in_service bool) ---

# Row 8
$ws.Cells.Item(8,3).Value = 14
$ws.Cells.Item(8,4).Value = 11

# Row 9
$ws.Cells.Item(9,3).Value = 16

# Row 10
$ws.Cells.Item(10,3).Value = 5
$ws.Cells.Item(10,4).Value = 12

# Row 11
$ws.Cells.Item(11,3).Value = 5
$ws.Cells.Item(11,4).Value = 9
$ws.Cells.Item(11,5).Value = $true

# Row 12
$ws.Cells.Item(12,3).Value = 10

# Row 13
$ws.Cells.Item(13,4).Value = 8

# Row 14
$ws.Cells.Item(14,3).Value = 9
$ws.Cells.Item(14,4).Value = 11
$ws.Cells.Item(14,5).Value = $true

# Row 15
$ws.Cells.Item(15,3).Value = 7
$ws.Cells.Item(15,4).Value = 11
$ws.Cells.Item(15,5).Value = $false

# --- Add two new contingency rows (line7, line8) ---

# Row 16
$ws.Cells.Item(15,1).Copy()
$ws.Cells.Item(16,1).PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "line7"
$ws.Cells.Item(16,3).Value = 5
$ws.Cells.Item(16,4).Value = 7
$ws.Cells.Item(16,5).Value = $true

# Row 17
$ws.Cells.Item(15,1).Copy()
$ws.Cells.Item(17,1).PasteSpecial(-4122)
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "line8"
$ws.Cells.Item(17,3).Value = 8
$ws.Cells.Item(17,4).Value = 5
$ws.Cells.Item(17,5).Value = $true
